$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update label text in A2 (shared string "Fnc" -> "Fc2")
$ws.Range("A2").Value = "Fc2"

# Update numeric results for row 2 (B2:F2); G2 stays 0
$ws.Range("B2").Value = 3.7763
$ws.Range("C2").Value = 11707.8297
$ws.Range("D2").Value = 3787.6483
$ws.Range("E2").Value = 4276.1202
$ws.Range("F2").Value = 3903.5288
